# Applies the cryptos price/volume update described by the commit diff.
# Most "Price" (column D) values look numeric, so they are written with a
# leading apostrophe to force Excel to store them as text (matching the
# original inline-string cell type) instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.403.84'
$ws.Range("E2").Value = '  -3.92%  '
$ws.Range("D3").Value = '1.950.83'
$ws.Range("E3").Value = '  -2.65%  '
$ws.Range("D4").Value = '''1.012'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''320.60'
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("D6").Value = '''1.010'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '''0.4754'
$ws.Range("E7").Value = '  -5.15%  '
$ws.Range("D8").Value = '''0.4069'
$ws.Range("E8").Value = '  -3.66%  '
$ws.Range("D9").Value = '''53.48'
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").Value = '''0.08469'
$ws.Range("E10").Value = '  -6.17%  '
$ws.Range("D11").Value = '''1.054'
$ws.Range("E11").Value = '  -5.73%  '
$ws.Range("D12").Value = '''22.01'
$ws.Range("E12").Value = '  -5.50%  '
$ws.Range("D13").Value = '1.980.45'
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("D14").Value = '''7.587'
$ws.Range("E14").Value = '  -5.63%  '
$ws.Range("D15").Value = '''6.155'
$ws.Range("E15").Value = '  -4.93%  '
$ws.Range("D16").Value = '''1.013'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '''0.00001073'
$ws.Range("E17").Value = '  -3.79%  '
$ws.Range("D18").Value = '''89.04'
$ws.Range("E18").Value = '  -5.81%  '
$ws.Range("D19").Value = '''0.06615'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").Value = '''18.67'
$ws.Range("E20").Value = '  -5.13%  '
$ws.Range("D21").Value = '''1.011'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '''5.811'
$ws.Range("E22").Value = '  -2.69%  '
$ws.Range("D23").Value = '28.415.96'
$ws.Range("E23").Value = '  -3.97%  '
$ws.Range("E24").Value = '  -4.14%  '
$ws.Range("D25").Value = '''2.290'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").Value = '2.198.14'
$ws.Range("E26").Value = '  -2.70%  '
$ws.Range("D27").Value = '''153.96'
$ws.Range("E27").Value = '  -3.17%  '
$ws.Range("D28").Value = '''20.16'
$ws.Range("E28").Value = '  -2.79%  '
$ws.Range("D29").Value = '''5.952'
$ws.Range("E29").Value = '  -6.40%  '
$ws.Range("D30").Value = '''2.150'
$ws.Range("E30").Value = '  -6.53%  '
$ws.Range("D31").Value = '''123.59'
$ws.Range("E31").Value = '  -3.67%  '
$ws.Range("D32").Value = '''0.9851'
$ws.Range("E32").Value = '  -6.95%  '
$ws.Range("D33").Value = '''0.09568'
$ws.Range("E33").Value = '  -4.04%  '
$ws.Range("D34").Value = '''1.445'
$ws.Range("E34").Value = '  -7.88%  '
$ws.Range("D35").Value = '''5.594'
$ws.Range("E35").Value = '  -4.17%  '
$ws.Range("D36").Value = '''3.655'
$ws.Range("E36").Value = '  -3.74%  '
$ws.Range("E37").Value = '  -5.46%  '
$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").Value = '''8.802'
$ws.Range("E38").Value = '  -5.29%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.06203'
$ws.Range("E39").Value = '  -3.53%  '
$ws.Range("D40").Value = '''1.253'
$ws.Range("E40").Value = '  -3.94%  '
$ws.Range("D41").Value = '''0.6217'
$ws.Range("E41").Value = '  -5.15%  '
$ws.Range("E42").Value = '  -5.08%  '
$ws.Range("D43").Value = '''1.010'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '''0.1916'
$ws.Range("E44").Value = '  -6.47%  '
$ws.Range("D45").Value = '''1.333'
$ws.Range("E45").Value = '  +1.75%  '
$ws.Range("D46").Value = '''0.5954'
$ws.Range("E46").Value = '  -6.44%  '
$ws.Range("D47").Value = '''12.94'
$ws.Range("E47").Value = '  -4.28%  '
$ws.Range("D48").Value = '''2.054'
$ws.Range("E48").Value = '  -6.60%  '
$ws.Range("D49").Value = '''3.403'
$ws.Range("E49").Value = '  -3.15%  '
$ws.Range("D50").Value = '''0.00000000328'
$ws.Range("E50").Value = '  -3.03%  '
$ws.Range("D51").Value = '''0.06806'
$ws.Range("E51").Value = '  -2.68%  '
